$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("1998D000002", "Jubiter Gloria", "Jubiter Michael, Mcgee Michael", "1998-02-01"),
    @("1998D000003", "Berryhill Cloris", "Berryhill Kenneth E", "1998-02-01"),
    @("1998D000005", "Leeper Stacey", "Leeper Glenn", "1998-02-01"),
    @("1998D000006", "Allen Aubrey", "Ross Barbara", "1998-02-01"),
    @("1998D000007", "Cartier Diana", "Delon Michael", "1998-02-01"),
    @("1998D000008", "Koons Sandra", "Koons Sam", "1998-02-01")
)

$row = 3
foreach ($rec in $data) {
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]

    # Force the filing-date column to be stored as text (matching the
    # source data, which keeps dates as plain "yyyy-mm-dd" strings)
    # instead of Excel auto-converting it to a date serial number.
    $ws.Cells.Item($row, 4).NumberFormat = "@"
    $ws.Cells.Item($row, 4).Value = $rec[3]

    $row++
}
